$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 200, pushing existing rows 200+ down to 202+.
$ws.Rows("200:201").Insert()

# Populate the two new rows with the new record (identical data in both rows).
$rows = @(200, 201)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 1).Value = 5
    $ws.Cells.Item($r, 2).Value = "Macroferia Regional de Talca"
    $ws.Cells.Item($r, 3).Value = "Maule"
    $ws.Cells.Item($r, 4).Value = 45126
    $ws.Cells.Item($r, 5).Value = 7
    $ws.Cells.Item($r, 6).Value = 100112017
    $ws.Cells.Item($r, 7).Value = "Apio"
    $ws.Cells.Item($r, 8).Value = "Americana (o)"
    $ws.Cells.Item($r, 9).Value = "Primera"
    $ws.Cells.Item($r, 10).Value = 700
    $ws.Cells.Item($r, 11).Value = 6000
    $ws.Cells.Item($r, 12).Value = 6000
    $ws.Cells.Item($r, 13).Value = 6000
    $ws.Cells.Item($r, 14).Value = "`$/docena de matas"
    $ws.Cells.Item($r, 15).Value = "Provincia del Elquí"
    $ws.Cells.Item($r, 16).Value = 1000
    $ws.Cells.Item($r, 17).Value = 6
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}
